$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1315
$ws.Range("I11").Value = 1315
$ws.Range("K11").Value = 1315
$ws.Range("M11").Value = -1175
$ws.Range("H17").Value = 1542.5555
$ws.Range("J17").Value = 1573.96
$ws.Range("L17").Value = 4721.88
$ws.Range("N17").Value = -5057.88
$ws.Range("H70").Value = 2807.9167
$ws.Range("I70").Value = 1799.6
$ws.Range("K70").Value = 5398.799999999999
$ws.Range("M70").Value = -5128.799999999999
$ws.Range("H73").Value = 2807.9167
$ws.Range("I73").Value = 1799.6
$ws.Range("K73").Value = 5398.799999999999
$ws.Range("M73").Value = -4462.799999999999
$ws.Range("H80").Value = 2334.75
$ws.Range("I80").Value = 938.5
$ws.Range("J80").Value = 2800.1667
$ws.Range("K80").Value = 2815.5
$ws.Range("L80").Value = 8400.500100000001
$ws.Range("M80").Value = -1817.5
$ws.Range("N80").Value = -10396.5001
$ws.Range("H83").Value = 2334.75
$ws.Range("I83").Value = 938.5
$ws.Range("J83").Value = 2800.1667
$ws.Range("K83").Value = 8446.5
$ws.Range("L83").Value = 25201.5003
$ws.Range("M83").Value = -3454.5
$ws.Range("N83").Value = -35185.5003
$ws.Range("H86").Value = 1791.5143
$ws.Range("I86").Value = 2135.4443
$ws.Range("K86").Value = 2135.4443
$ws.Range("M86").Value = -1012.4443
$ws.Range("H89").Value = 1791.5143
$ws.Range("I89").Value = 2135.4443
$ws.Range("K89").Value = 10677.2215
$ws.Range("M89").Value = -5061.2215
$ws.Range("H125").Value = 7410013.5
$ws.Range("I125").Value = 1449.5
$ws.Range("K125").Value = 13045.5
$ws.Range("M125").Value = -10585.5
$ws.Range("H138").Value = 3134.6785
$ws.Range("I138").Value = 2656
$ws.Range("J138").Value = 3247.3088
$ws.Range("K138").Value = 7968
$ws.Range("L138").Value = 9741.9264
$ws.Range("M138").Value = -2828
$ws.Range("N138").Value = -20021.9264

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3704556
$ws.Range("I2").Value = 5556264
$ws.Range("J2").Value = 1140
$ws.Range("K2").Value = 5556264
$ws.Range("L2").Value = 1140
$ws.Range("M2").Value = -5556151
$ws.Range("N2").Value = -1366
$ws.Range("H32").Value = 4890.552
$ws.Range("J32").Value = 8258
$ws.Range("L32").Value = 8258
$ws.Range("N32").Value = -8832
$ws.Range("H45").Value = 7574434
$ws.Range("I45").Value = 11067326
$ws.Range("J45").Value = 6500.6665
$ws.Range("K45").Value = 11067326
$ws.Range("L45").Value = 6500.6665
$ws.Range("M45").Value = -11066949
$ws.Range("N45").Value = -7254.6665
$ws.Range("H61").Value = 3895.8
$ws.Range("I61").Value = 3934.0715
$ws.Range("K61").Value = 3934.0715
$ws.Range("M61").Value = -3722.0715
$ws.Range("H74").Value = 59992.043
$ws.Range("I74").Value = 4558.7026
$ws.Range("J74").Value = 265095.4
$ws.Range("K74").Value = 4558.7026
$ws.Range("L74").Value = 265095.4
$ws.Range("M74").Value = -3684.7026
$ws.Range("N74").Value = -266843.4
$ws.Range("H77").Value = 59992.043
$ws.Range("I77").Value = 4558.7026
$ws.Range("J77").Value = 265095.4
$ws.Range("K77").Value = 22793.513
$ws.Range("L77").Value = 1325477
$ws.Range("M77").Value = -18425.513
$ws.Range("N77").Value = -1334213
$ws.Range("H116").Value = 3704556
$ws.Range("I116").Value = 5556264
$ws.Range("J116").Value = 1140
$ws.Range("K116").Value = 5556264
$ws.Range("L116").Value = 1140
$ws.Range("M116").Value = -5553970
$ws.Range("N116").Value = -5728
$ws.Range("H122").Value = 994989.5
$ws.Range("I122").Value = 2923.125
$ws.Range("J122").Value = 4169602
$ws.Range("K122").Value = 8769.375
$ws.Range("L122").Value = 12508806
$ws.Range("M122").Value = -6319.375
$ws.Range("N122").Value = -12513706
$ws.Range("H136").Value = 3895.8
$ws.Range("I136").Value = 3934.0715
$ws.Range("K136").Value = 11802.2145
$ws.Range("M136").Value = -9252.2145

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3704556
$ws.Range("I3").Value = 5556264
$ws.Range("J3").Value = 1140
$ws.Range("K3").Value = 5556264
$ws.Range("L3").Value = 1140
$ws.Range("M3").Value = -5556150
$ws.Range("N3").Value = -1368
$ws.Range("H74").Value = 22449.5
$ws.Range("I74").Value = 20999
$ws.Range("J74").Value = 23900
$ws.Range("K74").Value = 20999
$ws.Range("L74").Value = 23900
$ws.Range("M74").Value = -20063
$ws.Range("N74").Value = -25772
$ws.Range("H77").Value = 22449.5
$ws.Range("I77").Value = 20999
$ws.Range("J77").Value = 23900
$ws.Range("K77").Value = 62997
$ws.Range("L77").Value = 71700
$ws.Range("M77").Value = -58317
$ws.Range("N77").Value = -81060
$ws.Range("H92").Value = 44999.5
$ws.Range("J92").Value = 44999.5
$ws.Range("L92").Value = 44999.5
$ws.Range("N92").Value = -49991.5
$ws.Range("H99").Value = 11068824
$ws.Range("I99").Value = 15986292
$ws.Range("K99").Value = 15986292
$ws.Range("M99").Value = -15984794
$ws.Range("H107").Value = 3402341.2
$ws.Range("I107").Value = 4202643
$ws.Range("K107").Value = 4202643
$ws.Range("M107").Value = -4200723

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21941.467
$ws.Range("I31").Value = 1285.1724
$ws.Range("K31").Value = 1285.1724
$ws.Range("M31").Value = -990.1723999999999
$ws.Range("H34").Value = 21941.467
$ws.Range("I34").Value = 1285.1724
$ws.Range("K34").Value = 1285.1724
$ws.Range("M34").Value = -1083.1724
$ws.Range("H69").Value = 45049.25
$ws.Range("I69").Value = 17249
$ws.Range("K69").Value = 17249
$ws.Range("M69").Value = -16500
$ws.Range("H72").Value = 45049.25
$ws.Range("I72").Value = 17249
$ws.Range("K72").Value = 51747
$ws.Range("M72").Value = -48003
$ws.Range("H133").Value = 40326
$ws.Range("J133").Value = 40326
$ws.Range("L133").Value = 40326
$ws.Range("N133").Value = -45386

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 794584.6
$ws.Range("I97").Value = 953301.1
$ws.Range("J97").Value = 1002.2
$ws.Range("K97").Value = 953301.1
$ws.Range("L97").Value = 1002.2
$ws.Range("M97").Value = -952805.1
$ws.Range("N97").Value = -1994.2
$ws.Range("H102").Value = 3675395.8
$ws.Range("I102").Value = 4445767
$ws.Range("J102").Value = 2070455.6
$ws.Range("K102").Value = 4445767
$ws.Range("L102").Value = 2070455.6
$ws.Range("M102").Value = -4444145
$ws.Range("N102").Value = -2073699.6
$ws.Range("H107").Value = 1197.1111
$ws.Range("I107").Value = 1574
$ws.Range("K107").Value = 1574
$ws.Range("M107").Value = 346
$ws.Range("H126").Value = 9909512
$ws.Range("I126").Value = 4548083
$ws.Range("J126").Value = 27780944
$ws.Range("K126").Value = 13644249
$ws.Range("L126").Value = 83342832
$ws.Range("M126").Value = -13641779
$ws.Range("N126").Value = -83347772
$ws.Range("H132").Value = 3266.1072
$ws.Range("I132").Value = 2863.5386
$ws.Range("J132").Value = 8499.5
$ws.Range("K132").Value = 8590.6158
$ws.Range("L132").Value = 25498.5
$ws.Range("M132").Value = -6060.6158
$ws.Range("N132").Value = -30558.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3560.6667
$ws.Range("J68").Value = 3600
$ws.Range("L68").Value = 3600
$ws.Range("N68").Value = -5098
$ws.Range("H71").Value = 3560.6667
$ws.Range("J71").Value = 3600
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -25488
$ws.Range("H82").Value = 4632132.5
$ws.Range("I82").Value = 5558239.5
$ws.Range("J82").Value = 1598.5
$ws.Range("K82").Value = 5558239.5
$ws.Range("L82").Value = 1598.5
$ws.Range("M82").Value = -5557878.5
$ws.Range("N82").Value = -2320.5
$ws.Range("H85").Value = 4632132.5
$ws.Range("I85").Value = 5558239.5
$ws.Range("J85").Value = 1598.5
$ws.Range("K85").Value = 5558239.5
$ws.Range("L85").Value = 1598.5
$ws.Range("M85").Value = -5556991.5
$ws.Range("N85").Value = -4094.5
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H132").Value = 4002.25
$ws.Range("I132").Value = 3212.7727
$ws.Range("K132").Value = 9638.3181
$ws.Range("M132").Value = -7108.3181
$ws.Range("H136").Value = 46596.28
$ws.Range("I136").Value = 59945.23
$ws.Range("K136").Value = 179835.69
$ws.Range("M136").Value = -177285.69

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 83335840
$ws.Range("I81").Value = 83335840
$ws.Range("K81").Value = 166671680
$ws.Range("M81").Value = -166670619
$ws.Range("H84").Value = 83335840
$ws.Range("I84").Value = 83335840
$ws.Range("K84").Value = 833358400
$ws.Range("M84").Value = -833353096
$ws.Range("H104").Value = 41000
$ws.Range("J104").Value = 41000
$ws.Range("L104").Value = 41000
$ws.Range("N104").Value = -47988
$ws.Range("H113").Value = 721.3929000000001
$ws.Range("I113").Value = 688.5333000000001
$ws.Range("J113").Value = 759.3077
$ws.Range("K113").Value = 2065.5999
$ws.Range("L113").Value = 2277.9231
$ws.Range("M113").Value = 104.4000999999998
$ws.Range("N113").Value = -6617.9231
$ws.Range("H128").Value = 56536.25
$ws.Range("J128").Value = 60715
$ws.Range("L128").Value = 60715
$ws.Range("N128").Value = -70675
$ws.Range("H129").Value = 39563
$ws.Range("J129").Value = 39844.5
$ws.Range("L129").Value = 39844.5
$ws.Range("N129").Value = -49844.5
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

Write-Output "Applied all Hyperion Profits market-data updates."